$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56 (old rows 56-69 shift down to 57-70)
$ws.Rows.Item(56).Insert()

$ws.Cells.Item(56,1).Value  = 10
$ws.Cells.Item(56,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(56,3).Value  = "La Araucanía"
$ws.Cells.Item(56,4).Value  = 44634
$ws.Cells.Item(56,5).Value  = 9
$ws.Cells.Item(56,6).Value  = "Fruta"
$ws.Cells.Item(56,7).Value  = 100108
$ws.Cells.Item(56,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(56,9).Value  = 100108003
$ws.Cells.Item(56,10).Value = "Maracuyá"
$ws.Cells.Item(56,11).Value = "Sin especificar"
$ws.Cells.Item(56,12).Value = "Primera"
$ws.Cells.Item(56,13).Value = 30
$ws.Cells.Item(56,14).Value = 45000
$ws.Cells.Item(56,15).Value = 45000
$ws.Cells.Item(56,16).Value = 45000
$ws.Cells.Item(56,17).Value = "$/caja 18 kilos"
$ws.Cells.Item(56,18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(56,19).Value = 2500
$ws.Cells.Item(56,20).Value = 18

# Insert a second new row at position 65 (what was old row 64, now at 65, shifts down to 66)
$ws.Rows.Item(65).Insert()

$ws.Cells.Item(65,1).Value  = 10
$ws.Cells.Item(65,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(65,3).Value  = "La Araucanía"
$ws.Cells.Item(65,4).Value  = 44385
$ws.Cells.Item(65,5).Value  = 9
$ws.Cells.Item(65,6).Value  = "Fruta"
$ws.Cells.Item(65,7).Value  = 100108
$ws.Cells.Item(65,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(65,9).Value  = 100108003
$ws.Cells.Item(65,10).Value = "Maracuyá"
$ws.Cells.Item(65,11).Value = "Sin especificar"
$ws.Cells.Item(65,12).Value = "Primera"
$ws.Cells.Item(65,13).Value = 200
$ws.Cells.Item(65,14).Value = 1200
$ws.Cells.Item(65,15).Value = 1200
$ws.Cells.Item(65,16).Value = 1200
$ws.Cells.Item(65,17).Value = "$/kilo"
$ws.Cells.Item(65,18).Value = "Perú"
$ws.Cells.Item(65,19).Value = 1200
$ws.Cells.Item(65,20).Value = 1
